$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Programs" SQL query (B2) - adds a CASE expression for the Website column
$nl = "`r`n"
$newProgramsQuery = 'SELECT DISTINCT ' + $nl + `
    '    prg.program_name AS "Program",' + $nl + `
    '  CASE' + $nl + `
    '    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym' + $nl + `
    '        ELSE prg.program_link' + $nl + `
    '    END  AS "Website",' + $nl + `
    '    prg.focus_area AS "Focus Area",' + $nl + `
    '    prg.cancer_type AS "Cancer Type",' + $nl + `
    ' CASE ' + $nl + `
    '        WHEN prg.data_link IS NOT NULL THEN prg.website       ' + $nl + `
    '        ELSE prg.data_link' + $nl + `
    '    END AS "Data Location Details"' + $nl + `
    'FROM ' + $nl + `
    '    df_program prg' + $nl + `
    'WHERE ' + $nl + `
    "     prg.cancer_type LIKE '%Thyroid Cancer%'" + $nl + `
    'ORDER BY ' + $nl + `
    '    lower(prg.program_name) ASC' + $nl + `
    'LIMIT 100;'

$ws.Range("B2").Value = $newProgramsQuery

# Update the view: scroll back to the top of the sheet and select C3
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C3").Select()
